# Update MEG model run orders + run count, per commit "update: new MEG model runs"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Run count for the meg/hmm row (A26:D26) changes from 41 to 44
$ws.Range("C26").Value = 44

# Order strings for meg/hmm (row 26) and all meg/dynemo rows (rows 32-41)
$ws.Range("D26").Value = "[0,6,2,5,7,3,4,1]"

$ws.Range("D32").Value = "[0,3,6,1,4,2,7,5]"
$ws.Range("D33").Value = "[3,5,4,1,6,7,2,0]"
$ws.Range("D34").Value = "[6,2,7,0,4,3,5,1]"
$ws.Range("D35").Value = "[7,3,5,1,4,2,0,6]"
$ws.Range("D36").Value = "[5,7,1,0,6,2,3,4]"
$ws.Range("D37").Value = "[7,5,3,2,1,4,0,6]"
$ws.Range("D38").Value = "[3,5,1,2,7,6,0,4]"
$ws.Range("D39").Value = "[7,3,2,1,6,4,0,5]"
$ws.Range("D40").Value = "[3,7,2,4,6,0,1,5]"
$ws.Range("D41").Value = "[2,6,4,1,7,3,5,0]"
